$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 54, shifting existing rows 54-103 down to 57-106.
$ws.Rows.Item(54).Insert()
$ws.Rows.Item(54).Insert()
$ws.Rows.Item(54).Insert()

# New row 54: August Red / Primera
$ws.Range("A54").Value = 1
$ws.Range("B54").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C54").Value = "Arica y Parinacota"
$ws.Range("D54").Value = 44981
$ws.Range("D54").NumberFormat = $ws.Range("D57").NumberFormat
$ws.Range("E54").Value = 15
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100103
$ws.Range("H54").Value = "Frutos de hueso (carozo)"
$ws.Range("I54").Value = 100103006
$ws.Range("J54").Value = "Nectarín"
$ws.Range("K54").Value = "August Red"
$ws.Range("L54").Value = "Primera"
$ws.Range("M54").Value = 250
$ws.Range("N54").Value = 24000
$ws.Range("O54").Value = 25000
$ws.Range("P54").Value = 24400
$ws.Range("Q54").Value = "$/bandeja 18 kilos granel"
$ws.Range("R54").Value = "Región de O'Higgins"
$ws.Range("S54").Value = 1356
$ws.Range("T54").Value = 18

# New row 55: June Pearl / Primera
$ws.Range("A55").Value = 1
$ws.Range("B55").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C55").Value = "Arica y Parinacota"
$ws.Range("D55").Value = 44981
$ws.Range("D55").NumberFormat = $ws.Range("D57").NumberFormat
$ws.Range("E55").Value = 15
$ws.Range("F55").Value = "Fruta"
$ws.Range("G55").Value = 100103
$ws.Range("H55").Value = "Frutos de hueso (carozo)"
$ws.Range("I55").Value = 100103006
$ws.Range("J55").Value = "Nectarín"
$ws.Range("K55").Value = "June Pearl"
$ws.Range("L55").Value = "Primera"
$ws.Range("M55").Value = 300
$ws.Range("N55").Value = 24000
$ws.Range("O55").Value = 25000
$ws.Range("P55").Value = 24500
$ws.Range("Q55").Value = "$/bandeja 18 kilos granel"
$ws.Range("R55").Value = "Región de O'Higgins"
$ws.Range("S55").Value = 1361
$ws.Range("T55").Value = 18

# New row 56: Super Queen / Primera
$ws.Range("A56").Value = 1
$ws.Range("B56").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C56").Value = "Arica y Parinacota"
$ws.Range("D56").Value = 44981
$ws.Range("D56").NumberFormat = $ws.Range("D57").NumberFormat
$ws.Range("E56").Value = 15
$ws.Range("F56").Value = "Fruta"
$ws.Range("G56").Value = 100103
$ws.Range("H56").Value = "Frutos de hueso (carozo)"
$ws.Range("I56").Value = 100103006
$ws.Range("J56").Value = "Nectarín"
$ws.Range("K56").Value = "Super Queen"
$ws.Range("L56").Value = "Primera"
$ws.Range("M56").Value = 400
$ws.Range("N56").Value = 24000
$ws.Range("O56").Value = 25000
$ws.Range("P56").Value = 24500
$ws.Range("Q56").Value = "$/bandeja 18 kilos granel"
$ws.Range("R56").Value = "Región de O'Higgins"
$ws.Range("S56").Value = 1361
$ws.Range("T56").Value = 18
